$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.42
$ws.Range("G2").Value = 44
$ws.Range("H2").Value = 1.9
$ws.Range("I2").Value = 44
$ws.Range("J2").Value = 2.84
$ws.Range("K2").Value = 3.65
$ws.Range("V2").Value = 1.48
$ws.Range("W2").Value = 1.23
$ws.Range("AC2").Value = 42
$ws.Range("N3").Value = 1.1
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 1.23
$ws.Range("R3").Value = 1.21
$ws.Range("S3").Value = 1.56
$ws.Range("F4").Value = 2.12
$ws.Range("G4").Value = 2.3
$ws.Range("H4").Value = 3.15
$ws.Range("I4").Value = 3.55
$ws.Range("S4").Value = 2.72
$ws.Range("T4").Value = 1.63
$ws.Range("U4").Value = 2.32
$ws.Range("V4").Value = 1.39
$ws.Range("W4").Value = 1.76
$ws.Range("Y4").Value = 16
$ws.Range("Z4").Value = 90
$ws.Range("AA4").Value = 900
$ws.Range("AD4").Value = 15
$ws.Range("AE4").Value = 42
$ws.Range("AJ4").Value = 120
$ws.Range("AL4").Value = 80
$ws.Range("AM4").Value = 330
$ws.Range("AO4").Value = 32
$ws.Range("G5").Value = 2.42
$ws.Range("H5").Value = 3.45
$ws.Range("K5").Value = 3.6
$ws.Range("W5").Value = 1.71
$ws.Range("X5").Value = 14
$ws.Range("Y5").Value = 14.5
$ws.Range("AA5").Value = 90
$ws.Range("AB5").Value = 10
$ws.Range("AC5").Value = 8.6
$ws.Range("AD5").Value = 19
$ws.Range("AE5").Value = 60
$ws.Range("AG5").Value = 12.5
$ws.Range("AI5").Value = 75
$ws.Range("AK5").Value = 30
$ws.Range("AL5").Value = 50
$ws.Range("AM5").Value = 140
$ws.Range("AN5").Value = 27
$ws.Range("G6").Value = 1.28
$ws.Range("H6").Value = 15
$ws.Range("I6").Value = 15.5
$ws.Range("J6").Value = 6.2
$ws.Range("K6").Value = 7
$ws.Range("N6").Value = 5.2
$ws.Range("P6").Value = 2.42
$ws.Range("Q6").Value = 1.65
$ws.Range("S6").Value = 2.58
$ws.Range("U6").Value = 1.72
$ws.Range("W6").Value = 4.6
$ws.Range("X6").Value = 36
$ws.Range("AD6").Value = 130
$ws.Range("AH6").Value = 80
$ws.Range("AJ6").Value = 9.6
$ws.Range("AN6").Value = 4.6
$ws.Range("I7").Value = 6.2
$ws.Range("V7").Value = 1.19
$ws.Range("X7").Value = 12.5
$ws.Range("AB7").Value = 7.8
$ws.Range("AF7").Value = 11
$ws.Range("AL7").Value = 130
$ws.Range("F8").Value = 2.22
$ws.Range("G8").Value = 2.4
$ws.Range("J8").Value = 3.1
$ws.Range("K8").Value = 3.4
$ws.Range("N8").Value = 3
$ws.Range("W8").Value = 1.71
$ws.Range("Z8").Value = 65
$ws.Range("AE8").Value = 380
$ws.Range("AI8").Value = 380
$ws.Range("AL8").Value = 300
$ws.Range("AM8").Value = 580
$ws.Range("AO8").Value = 300
$ws.Range("J9").Value = 8.199999999999999
$ws.Range("L9").Value = 1.29
$ws.Range("O9").Value = 1.19
$ws.Range("Q9").Value = 1.6
$ws.Range("W9").Value = 6
$ws.Range("Z9").Value = 330
$ws.Range("AH9").Value = 70
$ws.Range("F10").Value = 1.3
$ws.Range("G10").Value = 1.31
$ws.Range("J10").Value = 5.9
$ws.Range("K10").Value = 6.4
$ws.Range("N10").Value = 4.3
$ws.Range("O10").Value = 1.24
$ws.Range("P10").Value = 2.14
$ws.Range("Q10").Value = 1.71
$ws.Range("R10").Value = 1.45
$ws.Range("T10").Value = 2.36
$ws.Range("W10").Value = 4.2
$ws.Range("Y10").Value = 980
$ws.Range("AD10").Value = 55
$ws.Range("AE10").Value = 380
$ws.Range("AF10").Value = 7.4
$ws.Range("AG10").Value = 12
$ws.Range("AH10").Value = 980
$ws.Range("AJ10").Value = 11
$ws.Range("AK10").Value = 16
$ws.Range("AL10").Value = 980
$ws.Range("AN10").Value = 6
$ws.Range("F11").Value = 1.87
$ws.Range("I11").Value = 5.3
$ws.Range("N11").Value = 3.15
$ws.Range("P11").Value = 1.75
$ws.Range("Q11").Value = 2.16
$ws.Range("R11").Value = 1.28
$ws.Range("AB11").Value = 8
$ws.Range("AD11").Value = 22
$ws.Range("AE11").Value = 80
$ws.Range("AF11").Value = 11
$ws.Range("AI11").Value = 95
$ws.Range("AK11").Value = 23
$ws.Range("AL11").Value = 980
$ws.Range("AM11").Value = 170
$ws.Range("AN11").Value = 16
$ws.Range("F12").Value = 1.97
$ws.Range("I12").Value = 4.9
$ws.Range("K12").Value = 3.55
$ws.Range("N12").Value = 3
$ws.Range("P12").Value = 1.68
$ws.Range("Q12").Value = 2.18
$ws.Range("R12").Value = 1.25
$ws.Range("T12").Value = 1.96
$ws.Range("X12").Value = 12
$ws.Range("AB12").Value = 8.199999999999999
$ws.Range("AJ12").Value = 25
$ws.Range("AN12").Value = 20
$ws.Range("G13").Value = 1.1
$ws.Range("H13").Value = 34
$ws.Range("I13").Value = 46
$ws.Range("J13").Value = 15.5
$ws.Range("K13").Value = 19
$ws.Range("Q13").Value = 1.31
$ws.Range("R13").Value = 2.08
$ws.Range("T13").Value = 2.64
$ws.Range("U13").Value = 1.49
$ws.Range("W13").Value = 11
$ws.Range("Y13").Value = 150
$ws.Range("AD13").Value = 150
$ws.Range("AF13").Value = 9.199999999999999
$ws.Range("AK13").Value = 17.5
$ws.Range("AM13").Value = 520
$ws.Range("F14").Value = 3.75
$ws.Range("G14").Value = 4.2
$ws.Range("H14").Value = 2.02
$ws.Range("I14").Value = 2.2
$ws.Range("K14").Value = 3.95
$ws.Range("P14").Value = 1.65

Write-Output "Applied 159 cell updates"